$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated "Size (links)" counts in column C ---
$ws.Range("C4").Value2 = 9
$ws.Range("C5").Value2 = 12
$ws.Range("C6").Value2 = 40
$ws.Range("C13").Value2 = 23
$ws.Range("C16").Value2 = 17
$ws.Range("C19").Value2 = 5
$ws.Range("C20").Value2 = 39
$ws.Range("C21").Value2 = 11
$ws.Range("C25").Value2 = 36
$ws.Range("C26").Value2 = 16
$ws.Range("C30").Value2 = 27
$ws.Range("C32").Value2 = 9
$ws.Range("C34").Value2 = 22
$ws.Range("C37").Value2 = 8
$ws.Range("C42").Value2 = 27
$ws.Range("C43").Value2 = 45
$ws.Range("C45").Value2 = 1
$ws.Range("C46").Value2 = 15
$ws.Range("C53").Value2 = 22
$ws.Range("C54").Value2 = 8
$ws.Range("C56").Value2 = 3
$ws.Range("C67").Value2 = 38
$ws.Range("C70").Value2 = 19

# --- Remove the stale "Error Link Handling" header cell and its note ---
$ws.Range("G1").Clear()
$ws.Range("H44").Clear()

# --- Mark column E ("History Status v2.0") cell styles for every data row ---
$ws.Range("E2").Style = "Neutral"
$ws.Range("E3:E65").Style = "Bad"
$ws.Range("E67:E70").Style = "Bad"

# --- Update the active selection to reflect where work left off ---
$ws.Range("J8").Select()
